$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (month shift: Oct->Nov, Nov->Dec)
# Use a leading quote so Excel treats these "Month Year" strings as text
# instead of auto-converting them into date serial numbers, then clear
# the resulting cell style back to Normal so no formatting is introduced.
$ws.Range("A1").Value = "'November 2024"
$ws.Range("A1").Style = "Normal"
$ws.Range("G1").Value = "'December 2024"
$ws.Range("G1").Style = "Normal"

# Update the numeric data row
$ws.Range("A2").Value = 0.536
$ws.Range("B2").Value = 0.316
$ws.Range("C2").Value = -0.008
$ws.Range("D2").Value = -0.113
$ws.Range("E2").Value = -0.006
$ws.Range("F2").Value = 0.267
$ws.Range("G2").Value = 0.984
